$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after PA5 (row 6) so PA5 can be split into PA5a / PA5b,
# pushing PA6..SA5 down by one row (rows 7-17 -> 8-18).
$ws.Rows("7:7").Insert()

# Split the old "PA5 / October 09, 2024" entry into two rows.
$ws.Range("B6").Value = "October 11, 2024"
$ws.Range("A6").Value = "PA5a"
$ws.Range("A7").Value = "PA5b"
$ws.Range("B7").Value = "October 16, 2024"

# The old "PA8 / November 13, 2024" row (now at row 10 after the insert)
# is removed entirely, shifting everything below it back up by one row.
$ws.Rows("10:10").Delete()

# The assignments that used to be PA9/PA10/PA11 shift down one label each
# (PA9 -> PA8, PA10 -> PA9, PA11 -> PA10) now that PA8 is gone.
$ws.Range("A10").Value = "PA8"
$ws.Range("A11").Value = "PA9"
$ws.Range("A12").Value = "PA10"

# Update the active selection to match the saved view state.
$ws.Range("B10").Select()
